# Documentary change: missing row (2 Jul 2013 effort entry) was added to
# the "Effort R 1.0" sheet (first/active sheet) of the workbook.
#
# New row 24:
#   A24 = 41457 (date 02/07/2013, inherits the date style from column A)
#   B24 = 2.5   (Effort [h])
#   D24 = "Implementation tc14, variants tried" (Task)
# C24 is left empty (no "Additional Effort [h]" logged for this entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 41457
$ws.Range("B24").Value = 2.5
$ws.Range("D24").Value = "Implementation tc14, variants tried"

# Move/keep the selection on C24, matching the saved view state.
$null = $ws.Range("C24").Select()
